# BOM update (v4) - refresh connector rows for the new sensor/IO breakout
# layout: add 5V/GND/I2C/SPI/UART header rows, move the IMU (LSM9DS1) next
# to its new connector, swap the old 9-pin header + loose GND/VBECOUT rows
# for the ADS1115 ADC, rename RPI02 -> RPI2, and drop the MCP3004 line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final BOM contents, row by row (row 1 is the header row).
$data = @(
    @("Comment","Designator","Footprint","LCSC Part Number"),
    @("Conn_01x03_Socket","5V1","PinHeader_1x03_P2.54mm_Vertical","C541850"),
    @("100nF","C1,C2,C3","C_0603_1608Metric_Pad1.08x0.95mm_HandSolder","C66501"),
    @("C_Polarized","C4","CP_Elec_4x5.4","C72483"),
    @("10nF","C5","C_0603_1608Metric_Pad1.08x0.95mm_HandSolder","C519406"),
    @("LEDG","D1","LED_0603_1608Metric","C965806"),
    @("LEDY","D2","LED_0603_1608Metric","C965802"),
    @("LEDR","D3","LED_0603_1608Metric","C965798"),
    @("LEDI","D4,D5,D6,D7","LED_0603_1608Metric","C412284"),
    @("Conn_01x03_Socket","GND1","PinHeader_1x03_P2.54mm_Vertical","C541850"),
    @("Conn_01x04_Socket","I2C1","PinHeader_1x04_P2.54mm_Vertical","C2718488"),
    @("LSM9DS1","IMU1","LGA-24L_3x3.5mm_P0.43mm","C2655096"),
    @("Conn_01x03_Socket","J1,J3,J4,J7,J37","PinHeader_1x03_P2.54mm_Vertical","C541850"),
    @("Conn_01x04_Socket","J2","PinSocket_1x04_P2.54mm_Vertical","C2718488"),
    @("Conn_01x02_Pin","J3_F2_BAT1","AMASS_XT60-F_1x02_P7.20mm_Vertical","C98733"),
    @("Conn_01x02_Pin","J4_M0_MOT0","AMASS_XT60-F_1x02_P7.20mm_Vertical","C98734"),
    @("Conn_01x02_Pin","J4_M1_MOT1","AMASS_XT60-F_1x02_P7.20mm_Vertical","C98734"),
    @("Conn_01x02_Pin","J4_M2_MOT2","AMASS_XT60-F_1x02_P7.20mm_Vertical","C98734"),
    @("Conn_01x02_Pin","J4_M3_MOT3","AMASS_XT60-F_1x02_P7.20mm_Vertical","C98734"),
    @("Conn_01x04_Pin","J9","JST_XH_B4B-XH-A_1x04_P2.50mm_Vertical","C18077835"),
    @("ADS1115IDGS","MAN1","TSSOP-10_3x3mm_P0.5mm","C37593"),
    @("130","R1,R2,R3,R7,R8,R9,R10,R11","R_0603_1608Metric","C22796"),
    @("49.9","R4,R5,R6","R_0603_1608Metric","C23185"),
    @("7K5","R14","R_0603_1608Metric","C728597"),
    @("2K5","R15","R_0603_1608Metric","C304065"),
    @("Raspberry_Pi_2_3","RPI2","Raspberry_Pi_Zero_Socketed_THT_FaceDown_MountingHoles","C2977589"),
    @("Conn_01x04_Socket","SPI1","PinHeader_1x04_P2.54mm_Vertical","C2718488"),
    @("Conn_01x04_Socket","UART1","PinHeader_1x04_P2.54mm_Vertical","C2718488")
)

$rowCount = $data.Length
$colCount = 4

# Write every cell. Column-A "Comment" entries like "130"/"49.9" read as
# pure numbers, but this BOM stores every cell as text (matching the
# original file), so Excel's auto-detect would otherwise silently turn
# them into doubles (losing the "130" formatting / introducing float
# noise on "49.9"). Collect those coordinates and force them to text via a
# transient "@" format, then drop the format straight back to the sheet's
# untouched default so no stray number formatting is left behind.
$numericLike = New-Object System.Collections.ArrayList

for ($r = 0; $r -lt $rowCount; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $colCount; $c++) {
        $val = $row[$c]
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
            [void]$numericLike.Add($cell)
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}

foreach ($cell in $numericLike) {
    $cell.ClearFormats()
}
